$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Header row: only the last column header text changes (DatePosted -> DaysPostedAgo)
$ws.Range("F1").Value = "DaysPostedAgo"

# Data rows 2-11: JobTitle, Company, Location, JobType, CompanySize
$data = @(
    @("Entry-Level Developer - required to work in office", "Vaco Technology", "Virginia Beach, VA", "Full Time", "N/A"),
    @("Work From Home Entry Level Data Entry Associate", "Level Up Entry", "Philadelphia, PA, PA", "Full-Time/Part-Time", "N/A"),
    @("Entry-level Healthcare Digital Technology Developer", "Cognizant Technology", "Dallas, TX", "Full Time", "N/A"),
    @("Entry Level Software Developer - Dev10 Technology Development Program - NY", "Genesis10", "New York, NY", "Full-Time/Part-Time", "N/A"),
    @("Entry Level Software Developer - Dev10 Technology Development Program - WI", "Genesis10", "Milwaukee, WI", "Full Time", "N/A"),
    @("Entry Level Software Developer - Dev10 Technology Development Program - MN", "Genesis10", "Minneapolis, MN", "Full Time", "N/A"),
    @("Entry-level EAS Digital Technology Developer", "Cognizant Technology", "Dallas, TX", "Full Time", "N/A"),
    @("Entry Level Software Developer - Dev10 Technology Development Program - .NET - TX", "Genesis10", "Dallas, TX", "Full Time", "N/A"),
    @("Entry Level Software Developer - Dev10 Technology Development Program - .NET - MN", "Genesis10", "Minneapolis, MN", "Full Time", "N/A"),
    @("Software Developer - Entry Level", "Revature", "Reston, VA", "Full Time", "N/A")
)

$r = 2
foreach ($row in $data) {
    $ws.Cells.Item($r, 1).Value = $row[0]
    $ws.Cells.Item($r, 2).Value = $row[1]
    $ws.Cells.Item($r, 3).Value = $row[2]
    $ws.Cells.Item($r, 4).Value = $row[3]
    $ws.Cells.Item($r, 5).Value = $row[4]
    $r = $r + 1
}

# Reflect the active selection shown in the saved file (F1 selected)
$ws.Range("F1").Select()
